$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-parsed by Excel as a number
# are written with NumberFormat forced to Text ("@") first, so they remain strings
# (matching the source inlineStr cells), same as the diff data.

$ws.Range("D2").Value = "60.569.82"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "2.595.94"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.66"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.50"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("E8").Value = "  +2.87%  "
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E10").Value = "  +1.97%  "
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.129"
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").Value = "3.052.27"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "60.519.14"
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.71"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").Value = "2.600.90"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.74"
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "351.34"
$ws.Range("E19").Value = "  +1.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.58"
$ws.Range("E20").Value = "  +2.46%  "
$ws.Range("E21").Value = "  +1.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "60.95"
$ws.Range("E23").Value = "  +1.04%  "
$ws.Range("E24").Value = "  +2.20%  "
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("D26").Value = "2.712.81"
$ws.Range("E26").Value = "  +0.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").Value = "0.0₃0842"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.29"
$ws.Range("E31").Value = "  +9.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.39"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("E33").Value = "  +2.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.95"
$ws.Range("E34").Value = "  -3.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.15"
$ws.Range("E35").Value = "  +3.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.20"
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.911"
$ws.Range("E37").Value = "  +6.36%  "
$ws.Range("E38").Value = "  +2.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.79"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("E40").Value = "  +1.49%  "
$ws.Range("E41").Value = "  -1.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "285.32"
$ws.Range("E42").Value = "  -4.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.101"
$ws.Range("E43").Value = "  +1.63%  "
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.56"
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("E48").Value = "  +1.04%  "
$ws.Range("E49").Value = "  -1.57%  "
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.14"
$ws.Range("E51").Value = "  +7.92%  "
